$wb = $excel.ActiveWorkbook

$wsBom = $wb.Worksheets.Item("BOM Report")
$wsInfo = $wb.Worksheets.Item("Project Information")

# --- BOM Report sheet ---

# Revision bumped from V1.0 to V1.1
$wsBom.Range("C5").Value = "V1.1"

# Manufacturer renamed: ON Semiconductor -> ONSEMI (D1 diode row)
$wsBom.Range("F18").Value = "ONSEMI"

# Connector (J3) swapped from GCT MICRO B USB to Molex USB2.0 Mini AB
$wsBom.Range("B20").Value = "CONN RCPT USB2.0 MINI AB 5P R/A"
$wsBom.Range("F20").Value = "Molex"
$wsBom.Range("G20").Value = "0565790519"
$wsBom.Range("J20").Value = "WM17121-ND"
$wsBom.Range("K20").Value = 2.04

# Designator list for the 0R jumper row now includes R10
$wsBom.Range("C24").Value = "R1, R5, R10"
$wsBom.Range("H24").Value = 3

# Distributor name/header fixes
$wsBom.Range("I35").Value = "Mouser"
$wsBom.Range("J9").Value = "02. distributor PN"

# --- Project Information sheet ---
# Leading apostrophe forces these number/date-looking values to remain
# plain text (matching the original quotePrefix-styled cells) instead of
# being auto-converted into numeric/date serials by Excel.
$wsInfo.Range("B7").Value = "'48"
$wsInfo.Range("B8").Value = "'1:15 PM"
$wsInfo.Range("B9").Value = "'10/13/2020"
$wsInfo.Range("B10").Value = "'10/13/2020 1:15 PM"

$wb.Application.Calculate()
